# Updated cryptos list - applies Price (D) and Volume(1h) (E) changes for rows 2-51
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.858.81"
$ws.Range("E2").Value = "  +3.19%  "

$ws.Range("D3").Value = "1.869.48"
$ws.Range("E3").Value = "  +2.76%  "

$ws.Range("E4").Value = "  +3.46%  "

$ws.Range("D5").Value = "'325.32"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.74%  "

$ws.Range("D6").Value = "'1.039"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.09%  "

$ws.Range("D7").Value = "'0.4435"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +3.04%  "

$ws.Range("E8").Value = "  +3.51%  "

$ws.Range("D9").Value = "'0.07485"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.22%  "

$ws.Range("D10").Value = "'0.8869"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.10%  "

$ws.Range("D11").Value = "'21.82"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.36%  "

$ws.Range("D12").Value = "1.876.18"
$ws.Range("E12").Value = "  -11.83%  "

$ws.Range("D13").Value = "'5.572"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.75%  "

$ws.Range("D14").Value = "'6.777"
$ws.Range("D14").Style = "Normal"

$ws.Range("D15").Value = "'0.07239"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +4.11%  "

$ws.Range("D16").Value = "'83.93"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.35%  "

$ws.Range("D17").Value = "'1.044"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +3.00%  "

$ws.Range("D18").Value = "'0.000009192"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.25%  "

$ws.Range("D19").Value = "'1.039"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.08%  "

$ws.Range("E20").Value = "  +2.47%  "

$ws.Range("D21").Value = "27.869.73"
$ws.Range("E21").Value = "  +3.10%  "

$ws.Range("D22").Value = "'5.335"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.54%  "

$ws.Range("E23").Value = "  +3.34%  "

$ws.Range("D24").Value = "'1.985"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +5.10%  "

$ws.Range("D25").Value = "'158.96"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.01%  "

$ws.Range("D27").Value = "'1.994"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +4.62%  "

$ws.Range("D28").Value = "'5.352"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.35%  "

$ws.Range("D29").Value = "'117.95"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.68%  "

$ws.Range("D30").Value = "'0.09104"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.63%  "

$ws.Range("D31").Value = "'0.7796"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +4.04%  "

$ws.Range("D32").Value = "'3.118"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +10.63%  "

$ws.Range("D33").Value = "'1.218"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.49%  "

$ws.Range("D34").Value = "'4.588"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.46%  "

$ws.Range("D35").Value = "'1.041"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.17%  "

$ws.Range("D36").Value = "'1.160"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.38%  "

$ws.Range("D37").Value = "'0.02001"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.75%  "

$ws.Range("D38").Value = "'0.05366"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.40%  "

$ws.Range("D39").Value = "'2.868"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +4.60%  "

$ws.Range("D40").Value = "'0.5211"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.92%  "

$ws.Range("E41").Value = "  +2.44%  "

$ws.Range("D42").Value = "'6.925"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +6.54%  "

$ws.Range("D43").Value = "'8.691"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.98%  "

$ws.Range("D44").Value = "'110.05"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.81%  "

$ws.Range("D45").Value = "'10.72"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.94%  "

$ws.Range("E46").Value = "  +5.10%  "

$ws.Range("E47").Value = "  +3.26%  "

$ws.Range("D48").Value = "'0.06473"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +4.02%  "

$ws.Range("D49").Value = "'1.905"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.77%  "

$ws.Range("D50").Value = "'40.03"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.97%  "

$ws.Range("D51").Value = "'64.89"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.96%  "
